$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 21:39"

# Row 4
$ws.Cells.Item(4, 2).Value = 6612515
$ws.Cells.Item(4, 3).Value = 24355
$ws.Cells.Item(4, 4).Value = 3886185
$ws.Cells.Item(4, 5).Value = 2529474
$ws.Cells.Item(4, 7).Value = 529
$ws.Cells.Item(4, 8).Value = 196856

# Row 5
$ws.Cells.Item(5, 2).Value = 4657379
$ws.Cells.Item(5, 3).Value = 97654
$ws.Cells.Item(5, 4).Value = 3621438
$ws.Cells.Item(5, 5).Value = 958435
$ws.Cells.Item(5, 7).Value = 1202
$ws.Cells.Item(5, 8).Value = 77506

# Row 12
$ws.Cells.Item(12, 2).Value = 576697
$ws.Cells.Item(12, 3).Value = 4708

# Row 24
$ws.Cells.Item(24, 2).Value = 259720
$ws.Cells.Item(24, 3).Value = 1613
$ws.Cells.Item(24, 5).Value = 16997

# Row 43
$ws.Cells.Item(43, 1).Value = "Marruecos"
$ws.Cells.Item(43, 2).Value = 82197
$ws.Cells.Item(43, 3).Value = 2430
$ws.Cells.Item(43, 4).Value = 64194
$ws.Cells.Item(43, 5).Value = 16479
$ws.Cells.Item(43, 7).Value = 33
$ws.Cells.Item(43, 8).Value = 1524

# Row 44
$ws.Cells.Item(44, 1).Value = "Guatemala"
$ws.Cells.Item(44, 2).Value = 81009
$ws.Cells.Item(44, 3).Value = 703
$ws.Cells.Item(44, 4).Value = 69703
$ws.Cells.Item(44, 5).Value = 8377
$ws.Cells.Item(44, 7).Value = 11
$ws.Cells.Item(44, 8).Value = 2929

# Row 45
$ws.Cells.Item(45, 1).Value = "Paises Bajos"
$ws.Cells.Item(45, 2).Value = 79781
$ws.Cells.Item(45, 3).Value = 1270
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 7).Value = 3
$ws.Cells.Item(45, 8).Value = 6252

# Row 93
$ws.Cells.Item(93, 1).Value = "Albania"
$ws.Cells.Item(93, 2).Value = 11021
$ws.Cells.Item(93, 3).Value = 161
$ws.Cells.Item(93, 4).Value = 6443
$ws.Cells.Item(93, 5).Value = 4251
$ws.Cells.Item(93, 7).Value = 3
$ws.Cells.Item(93, 8).Value = 327

# Row 94
$ws.Cells.Item(94, 1).Value = "Hungria"
$ws.Cells.Item(94, 2).Value = 10909
$ws.Cells.Item(94, 3).Value = 718
$ws.Cells.Item(94, 4).Value = 4014
$ws.Cells.Item(94, 5).Value = 6264
$ws.Cells.Item(94, 7).Value = 1
$ws.Cells.Item(94, 8).Value = 631

# Row 104
$ws.Cells.Item(104, 2).Value = 8457
$ws.Cells.Item(104, 3).Value = 28
$ws.Cells.Item(104, 4).Value = 6120
$ws.Cells.Item(104, 5).Value = 2121
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 216

# Row 108
$ws.Cells.Item(108, 1).Value = "Tunez"
$ws.Cells.Item(108, 2).Value = 6259
$ws.Cells.Item(108, 3).Value = 377
$ws.Cells.Item(108, 4).Value = 1956
$ws.Cells.Item(108, 5).Value = 4200
$ws.Cells.Item(108, 7).Value = 4
$ws.Cells.Item(108, 8).Value = 103

# Row 109
$ws.Cells.Item(109, 1).Value = "Montenegro"
$ws.Cells.Item(109, 2).Value = 6222
$ws.Cells.Item(109, 4).Value = 4393
$ws.Cells.Item(109, 5).Value = 1715
$ws.Cells.Item(109, 8).Value = 114

# Row 110
$ws.Cells.Item(110, 2).Value = 5669
$ws.Cells.Item(110, 3).Value = 14
$ws.Cells.Item(110, 4).Value = 3720
$ws.Cells.Item(110, 5).Value = 1772
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = 177

# Row 113
$ws.Cells.Item(113, 2).Value = 5025
$ws.Cells.Item(113, 3).Value = 31
$ws.Cells.Item(113, 4).Value = 4165
$ws.Cells.Item(113, 5).Value = 762

# Row 115
$ws.Cells.Item(115, 1).Value = "Congo"
$ws.Cells.Item(115, 2).Value = 4928
$ws.Cells.Item(115, 3).Value = 37
$ws.Cells.Item(115, 4).Value = 3887
$ws.Cells.Item(115, 5).Value = 953
$ws.Cells.Item(115, 7).Value = 5
$ws.Cells.Item(115, 8).Value = 88

# Row 116
$ws.Cells.Item(116, 1).Value = "Hong Kong"
$ws.Cells.Item(116, 2).Value = 4926
$ws.Cells.Item(116, 3).Value = 12
$ws.Cells.Item(116, 4).Value = 4597
$ws.Cells.Item(116, 5).Value = 230
$ws.Cells.Item(116, 8).Value = 99

# Row 117
$ws.Cells.Item(117, 1).Value = "Mozambique"
$ws.Cells.Item(117, 2).Value = 4918
$ws.Cells.Item(117, 3).Value = 86
$ws.Cells.Item(117, 4).Value = 2899
$ws.Cells.Item(117, 5).Value = 1988
$ws.Cells.Item(117, 8).Value = 31

# Row 135
$ws.Cells.Item(135, 1).Value = "Jordania"
$ws.Cells.Item(135, 2).Value = 2945
$ws.Cells.Item(135, 3).Value = 206
$ws.Cells.Item(135, 4).Value = 2084
$ws.Cells.Item(135, 5).Value = 840
$ws.Cells.Item(135, 7).Value = 1
$ws.Cells.Item(135, 8).Value = 21

# Row 136
$ws.Cells.Item(136, 1).Value = "Mali"
$ws.Cells.Item(136, 2).Value = 2912
$ws.Cells.Item(136, 3).Value = 3
$ws.Cells.Item(136, 4).Value = 2271
$ws.Cells.Item(136, 5).Value = 513
$ws.Cells.Item(136, 8).Value = 128

# Row 137
$ws.Cells.Item(137, 1).Value = "Aruba"
$ws.Cells.Item(137, 2).Value = 2819
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 1407
$ws.Cells.Item(137, 5).Value = 1396
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 16

# Row 138
$ws.Cells.Item(138, 1).Value = "Bahamas"
$ws.Cells.Item(138, 2).Value = 2814
$ws.Cells.Item(138, 3).Value = 93
$ws.Cells.Item(138, 4).Value = 1220
$ws.Cells.Item(138, 5).Value = 1529
$ws.Cells.Item(138, 8).Value = 65

# Row 139
$ws.Cells.Item(139, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(139, 2).Value = 2777
$ws.Cells.Item(139, 3).Value = 79
$ws.Cells.Item(139, 4).Value = 762
$ws.Cells.Item(139, 5).Value = 1970
$ws.Cells.Item(139, 7).Value = 2
$ws.Cells.Item(139, 8).Value = 45

# Row 157
$ws.Cells.Item(157, 2).Value = 1520
$ws.Cells.Item(157, 3).Value = 3
$ws.Cells.Item(157, 5).Value = 261

# Row 165
$ws.Cells.Item(165, 1).Value = "Republica del Chad"
$ws.Cells.Item(165, 2).Value = 1081
$ws.Cells.Item(165, 3).Value = 30
$ws.Cells.Item(165, 4).Value = 934
$ws.Cells.Item(165, 5).Value = 68
$ws.Cells.Item(165, 8).Value = 79

# Row 166
$ws.Cells.Item(166, 1).Value = "Vietnam"
$ws.Cells.Item(166, 2).Value = 1060
$ws.Cells.Item(166, 3).Value = 1
$ws.Cells.Item(166, 4).Value = 902
$ws.Cells.Item(166, 5).Value = 123
$ws.Cells.Item(166, 8).Value = 35
